$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5701
$ws.Range("I64").Value = 5701
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 5701
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -5453
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5701
$ws.Range("I67").Value = 5701
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 5701
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -4843
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 12505915
$ws.Range("I111").Value = 15631272
$ws.Range("J111").Value = 4484.5
$ws.Range("K111").Value = 46893816
$ws.Range("L111").Value = 13453.5
$ws.Range("M111").Value = -46890749
$ws.Range("N111").Value = -19587.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1845.25
$ws.Range("I129").Value = 1845.25
$ws.Range("K129").Value = 5535.75
$ws.Range("M129").Value = -535.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3450.4333
$ws.Range("I132").Value = 3446.9285
$ws.Range("K132").Value = 10340.7855
$ws.Range("M132").Value = -7810.7855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3235
$ws.Range("I137").Value = 3147
$ws.Range("K137").Value = 9441
$ws.Range("M137").Value = -6891

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 59999
$ws.Range("J140").Value = 59999
$ws.Range("L140").Value = 59999
$ws.Range("N140").Value = -70359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2407220.8
$ws.Range("I32").Value = 2607323.2
$ws.Range("K32").Value = 2607323.2
$ws.Range("M32").Value = -2607036.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1776.4706
$ws.Range("I63").Value = 1730.7142
$ws.Range("K63").Value = 1730.7142
$ws.Range("M63").Value = -1044.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1776.4706
$ws.Range("I66").Value = 1730.7142
$ws.Range("K66").Value = 8653.571
$ws.Range("M66").Value = -5221.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 5567031.5
$ws.Range("I97").Value = 1065.1666
$ws.Range("J97").Value = 9277676
$ws.Range("K97").Value = 1065.1666
$ws.Range("L97").Value = 9277676
$ws.Range("M97").Value = -569.1666
$ws.Range("N97").Value = -9278668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 82000
$ws.Range("J139").Value = 82000
$ws.Range("L139").Value = 82000
$ws.Range("N139").Value = -92280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7248367.5
$ws.Range("J20").Value = 2122
$ws.Range("L20").Value = 2122
$ws.Range("N20").Value = -2616

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 30305736
$ws.Range("I86").Value = 1834.8077
$ws.Range("J86").Value = 142863090
$ws.Range("K86").Value = 1834.8077
$ws.Range("L86").Value = 142863090
$ws.Range("M86").Value = -711.8077000000001
$ws.Range("N86").Value = -142865336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 30305736
$ws.Range("I89").Value = 1834.8077
$ws.Range("J89").Value = 142863090
$ws.Range("K89").Value = 9174.038500000001
$ws.Range("L89").Value = 714315450
$ws.Range("M89").Value = -3558.038500000001
$ws.Range("N89").Value = -714326682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 169166740
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 169166740
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 507500220
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -507500786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 62.285713
$ws.Range("I38").Value = 49
$ws.Range("J38").Value = 72.25
$ws.Range("K38").Value = 147
$ws.Range("L38").Value = 216.75
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = -910.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10699.2
$ws.Range("I39").Value = 501
$ws.Range("J39").Value = 15069.857
$ws.Range("K39").Value = 1503
$ws.Range("L39").Value = 45209.571
$ws.Range("M39").Value = -1209
$ws.Range("N39").Value = -45797.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 819.2
$ws.Range("I46").Value = 524.5
$ws.Range("J46").Value = 1998
$ws.Range("K46").Value = 1573.5
$ws.Range("L46").Value = 5994
$ws.Range("M46").Value = -1482.5
$ws.Range("N46").Value = -6176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 223333600
$ws.Range("I51").Value = 400
$ws.Range("K51").Value = 1200
$ws.Range("M51").Value = -740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5666.6665
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 8000
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 24000
$ws.Range("M57").Value = -2441
$ws.Range("N57").Value = -25118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4993.6
$ws.Range("I80").Value = 4994
$ws.Range("K80").Value = 4994
$ws.Range("M80").Value = -3996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4993.6
$ws.Range("I83").Value = 4994
$ws.Range("K83").Value = 24970
$ws.Range("M83").Value = -19978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1775.6757
$ws.Range("I97").Value = 1800.9048
$ws.Range("J97").Value = 1742.5625
$ws.Range("K97").Value = 1800.9048
$ws.Range("L97").Value = 1742.5625
$ws.Range("M97").Value = -1304.9048
$ws.Range("N97").Value = -2734.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4321.1665
$ws.Range("I102").Value = 4268.5
$ws.Range("K102").Value = 4268.5
$ws.Range("M102").Value = -2646.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 40763.742
$ws.Range("J122").Value = 4495.273
$ws.Range("L122").Value = 13485.819
$ws.Range("N122").Value = -18385.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2749.0715
$ws.Range("I126").Value = 2818.8
$ws.Range("K126").Value = 8456.400000000001
$ws.Range("M126").Value = -5986.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6008.5625
$ws.Range("I40").Value = 5663.4
$ws.Range("K40").Value = 5663.4
$ws.Range("M40").Value = -5527.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3174.889
$ws.Range("I82").Value = 2914.8
$ws.Range("K82").Value = 2914.8
$ws.Range("M82").Value = -2553.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3174.889
$ws.Range("I85").Value = 2914.8
$ws.Range("K85").Value = 2914.8
$ws.Range("M85").Value = -1666.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4544.769
$ws.Range("I122").Value = 3065.2307
$ws.Range("K122").Value = 9195.6921
$ws.Range("M122").Value = -6745.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 65256.832
$ws.Range("J141").Value = 65256.832
$ws.Range("L141").Value = 65256.832
$ws.Range("N141").Value = -75616.83199999999
